$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Minimo Ruptura" values for column C (rows 2-4 y 7-11), usando
# el mismo color de relleno (theme accent2, tint ~0.8 -> RGB FBE5D6)
$fillColor = 14083579

$ws.Range("C2").Value = 4500
$ws.Range("C2").Interior.Color = $fillColor

$ws.Range("C3").Value = 4500
$ws.Range("C3").Interior.Color = $fillColor

$ws.Range("C4").Value = 4500
$ws.Range("C4").Interior.Color = $fillColor

$ws.Range("C7").Value = 7000
$ws.Range("C7").Interior.Color = $fillColor

$ws.Range("C8").Value = 7000
$ws.Range("C8").Interior.Color = $fillColor

$ws.Range("C9").Value = 9000
$ws.Range("C9").Interior.Color = $fillColor

$ws.Range("C10").Value = 7000
$ws.Range("C10").Interior.Color = $fillColor

$ws.Range("C11").Value = 6500
$ws.Range("C11").Interior.Color = $fillColor

# Actualizar valores de ruptura minima (B5, B6)
$ws.Range("B5").Value = 90
$ws.Range("B6").Value = 90

# Actualizar la celda seleccionada
$ws.Range("C12").Select()
